# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (fund-holding detail) right before the
#    "总计" (total) sheet, matching the layout already used by the other
#    quarterly sheets (基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名).
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" detail sheet by cloning "2021-Q4" (same
#    header row / formatting / page setup) and overwriting its data.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$templateSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template has 3 data rows (rows 2-4); the new sheet only needs 2.
$newSheet.Range("A4").EntireRow.Delete()

# Row 2: 006792 / 鹏华香港美国互联网股票（LOF）美元现汇
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "006792"
$newSheet.Range("C2").Value = "鹏华香港美国互联网股票（LOF）美元现汇"
$newSheet.Range("D2").Value = "1.43"
$newSheet.Range("E2").Value = "83.72"
$newSheet.Range("F2").Value = "2.73"
$newSheet.Range("G2").Value = "0.0390"
$newSheet.Range("B2:G2").ClearFormats()
$newSheet.Range("H2").Value = 10

# Row 3: 160644 / 鹏华香港美国互联网股票（LOF）人民币
$newSheet.Range("B3:G3").NumberFormat = "@"
$newSheet.Range("B3").Value = "160644"
$newSheet.Range("C3").Value = "鹏华香港美国互联网股票（LOF）人民币"
$newSheet.Range("D3").Value = "1.43"
$newSheet.Range("E3").Value = "83.72"
$newSheet.Range("F3").Value = "2.73"
$newSheet.Range("G3").Value = "0.0390"
$newSheet.Range("B3:G3").ClearFormats()
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2. Prepend a new summary row to "总计" for the 2022-Q1 quarter, shifting
#    the existing rows down and renumbering the index column (A).
#    NOTE: re-resolve the sheet by name — inserting the new sheet shifted
#    "总计"'s position, so the old $totalSheet reference now points at the
#    wrong (positionally-indexed) sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("A2").EntireRow.Insert()

# Re-apply the index-column formatting (border/alignment) from the row
# below, which still carries the original style.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.08

# Renumber the remaining index column values (0,1,2,3 -> 1,2,3,4).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
